# Insert two new rows at the top of the data block (row 99), pushing the
# existing rows 99:139 down to 101:141, then populate the two new rows with
# the new week's data (Fecha 2021-09-16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 99 onward down by two rows.
$ws.Rows("99:100").Insert()

# New row 99 - "Primera"
$ws.Range("A99").Value = 11
$ws.Range("B99").Value = "Vega Monumental Concepción"
$ws.Range("C99").Value = "Bíobío"
$ws.Range("D99").Value = "2021-09-16"
$ws.Range("D99").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E99").Value = 8
$ws.Range("F99").Value = 100112008
$ws.Range("G99").Value = "Coliflor"
$ws.Range("H99").Value = "Sin especificar"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 700
$ws.Range("L99").Value = 800
$ws.Range("M99").Value = 750
$ws.Range("N99").Value = "$/unidad"
$ws.Range("O99").Value = "Región Metropolitana"
$ws.Range("P99").Value = 750
$ws.Range("Q99").Value = 1
$ws.Range("R99").Value = "Hortaliza"

# New row 100 - "Segunda"
$ws.Range("A100").Value = 11
$ws.Range("B100").Value = "Vega Monumental Concepción"
$ws.Range("C100").Value = "Bíobío"
$ws.Range("D100").Value = "2021-09-16"
$ws.Range("D100").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E100").Value = 8
$ws.Range("F100").Value = 100112008
$ws.Range("G100").Value = "Coliflor"
$ws.Range("H100").Value = "Sin especificar"
$ws.Range("I100").Value = "Segunda"
$ws.Range("J100").Value = 500
$ws.Range("K100").Value = 600
$ws.Range("L100").Value = 600
$ws.Range("M100").Value = 600
$ws.Range("N100").Value = "$/unidad"
$ws.Range("O100").Value = "Región Metropolitana"
$ws.Range("P100").Value = 600
$ws.Range("Q100").Value = 1
$ws.Range("R100").Value = "Hortaliza"
